# Yêu cầu bài tập.docx - "cập nhật phần tạo folder cho controler" edit
#
# 1) Several numbered list lines were split across two runs (the leading
#    digit in one run, the ". <description>" text in a second run, both
#    sharing identical run formatting). Word had silently merged those
#    pairs back into a single run/string - reproduce that by doing a
#    plain Find & Replace of the split text with the joined text.
# 2) An empty trailing paragraph right after item "10. ..." was removed.

$d = $word.ActiveDocument

function Merge-SplitLine($oldText, $newText) {
    $range = $d.Content
    $range.Find.ClearFormatting()
    $found = $range.Find.Execute($oldText, $true, $true, $false, $false, $false, `
                                  $true, 1, $false, $newText, 2)
}

Merge-SplitLine "2. Thông tin hợp tác quốc tế (HTQT, HTDN) gồm 10 bảng" "2. Thông tin hợp tác quốc tế (HTQT, HTDN) gồm 10 bảng"
Merge-SplitLine "3. Thông tin khoa học công nghệ (KHCN) gồm 10 bảng" "3. Thông tin khoa học công nghệ (KHCN) gồm 10 bảng"
Merge-SplitLine "4. Nhóm người học (NH) gồm 10 bảng" "4. Nhóm người học (NH) gồm 10 bảng"
Merge-SplitLine "9. Nhóm ngành đào tạo (NDT) gồm 5 bảng" "9. Nhóm ngành đào tạo (NDT) gồm 5 bảng"

# Remove the empty paragraph that sits right after "10. Thông tin tuyển sinh (TS) gồm 2 bảng"
foreach ($p in $d.Paragraphs) {
    $text = $p.Range.Text
    $trimmed = $text.Trim("`r", "`a", "`n", " ")
    if ($trimmed -eq "") {
        $prev = $p.Previous()
        if ($prev -ne $null -and $prev.Range.Text.Contains("10. Thông tin tuyển sinh")) {
            $p.Range.Delete()
            break
        }
    }
}
